# Updates cryptocurrency price/volume data in the worksheet to reflect
# the latest scrape performed by the GitHub Actions workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.856.57'
$ws.Range("E2").Value = '  +1.29%  '
$ws.Range("D3").Value = '2.104.09'
$ws.Range("E3").Value = '  +2.24%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.38'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.83%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +2.14%  '
$ws.Range("E10").Value = '  +2.21%  '
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("D12").Value = '2.414.65'
$ws.Range("E12").Value = '  +2.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.51'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.74%  '
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.35%  '
$ws.Range("D17").Value = '2.113.66'
$ws.Range("E17").Value = '  +2.98%  '
$ws.Range("D18").Value = '37.801.72'
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("D21").Value = '0.0₃0823'
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("E28").Value = '  +2.85%  '
$ws.Range("E29").Value = '  -3.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.53'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.99%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.83%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0623'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.58'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.21%  '
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.82%  '
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0963'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("D43").Value = '1.475.79'
$ws.Range("E43").Value = '  +1.19%  '
$ws.Range("E44").Value = '  +1.02%  '
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -11.26%  '
$ws.Range("E47").Value = '  +1.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '15.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("E50").Value = '  +3.17%  '
$ws.Range("D51").Value = '2.301.25'
$ws.Range("E51").Value = '  +2.40%  '
